$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2204.611
$ws.Range("I62").Value = 2127.5454
$ws.Range("J62").Value = 2325.7144
$ws.Range("K62").Value = 2127.5454
$ws.Range("L62").Value = 2325.7144
$ws.Range("M62").Value = -1503.5454
$ws.Range("N62").Value = -3573.7144
$ws.Range("H64").Value = 3723.077
$ws.Range("I64").Value = 3312.5
$ws.Range("J64").Value = 4380
$ws.Range("K64").Value = 3312.5
$ws.Range("L64").Value = 4380
$ws.Range("M64").Value = -3064.5
$ws.Range("N64").Value = -4876
$ws.Range("H65").Value = 2204.611
$ws.Range("I65").Value = 2127.5454
$ws.Range("J65").Value = 2325.7144
$ws.Range("K65").Value = 10637.727
$ws.Range("L65").Value = 11628.572
$ws.Range("M65").Value = -7517.726999999999
$ws.Range("N65").Value = -17868.572
$ws.Range("H67").Value = 3723.077
$ws.Range("I67").Value = 3312.5
$ws.Range("J67").Value = 4380
$ws.Range("K67").Value = 3312.5
$ws.Range("L67").Value = 4380
$ws.Range("M67").Value = -2454.5
$ws.Range("N67").Value = -6096

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8377.929
$ws.Range("J63").Value = 12913.143
$ws.Range("L63").Value = 12913.143
$ws.Range("N63").Value = -14285.143
$ws.Range("H66").Value = 8377.929
$ws.Range("J66").Value = 12913.143
$ws.Range("L66").Value = 64565.715
$ws.Range("N66").Value = -71429.715
$ws.Range("H80").Value = 200026510
$ws.Range("I80").Value = 28289
$ws.Range("J80").Value = 333358700
$ws.Range("K80").Value = 28289
$ws.Range("L80").Value = 333358700
$ws.Range("M80").Value = -27291
$ws.Range("N80").Value = -333360696
$ws.Range("H83").Value = 200026510
$ws.Range("I83").Value = 28289
$ws.Range("J83").Value = 333358700
$ws.Range("K83").Value = 84867
$ws.Range("L83").Value = 1000076100
$ws.Range("M83").Value = -79875
$ws.Range("N83").Value = -1000086084
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 8898
$ws.Range("I82").Value = 2051.75
$ws.Range("K82").Value = 2051.75
$ws.Range("M82").Value = -1668.75
$ws.Range("H85").Value = 8898
$ws.Range("I85").Value = 2051.75
$ws.Range("K85").Value = 2051.75
$ws.Range("M85").Value = -725.75
$ws.Range("H86").Value = 2231.6
$ws.Range("I86").Value = 2390
$ws.Range("J86").Value = 1598
$ws.Range("K86").Value = 2390
$ws.Range("L86").Value = 1598
$ws.Range("M86").Value = -1267
$ws.Range("N86").Value = -3844
$ws.Range("H89").Value = 2231.6
$ws.Range("I89").Value = 2390
$ws.Range("J89").Value = 1598
$ws.Range("K89").Value = 11950
$ws.Range("L89").Value = 7990
$ws.Range("M89").Value = -6334
$ws.Range("N89").Value = -19222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1497.25
$ws.Range("I16").Value = 1550.2858
$ws.Range("J16").Value = 1423
$ws.Range("K16").Value = 1550.2858
$ws.Range("L16").Value = 1423
$ws.Range("M16").Value = -1263.2858
$ws.Range("N16").Value = -1997
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H113").Value = 1497.25
$ws.Range("I113").Value = 1550.2858
$ws.Range("J113").Value = 1423
$ws.Range("K113").Value = 1550.2858
$ws.Range("L113").Value = 1423
$ws.Range("M113").Value = 619.7141999999999
$ws.Range("N113").Value = -5763

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.18421
$ws.Range("I2").Value = 15.5
$ws.Range("J2").Value = 50.75
$ws.Range("K2").Value = 93
$ws.Range("L2").Value = 304.5
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = -530.5
$ws.Range("H14").Value = 31568.1
$ws.Range("I14").Value = 31568.1
$ws.Range("K14").Value = 94704.29999999999
$ws.Range("M14").Value = -94531.29999999999
$ws.Range("H56").Value = 104029.9
$ws.Range("I56").Value = 104029.9
$ws.Range("K56").Value = 104029.9
$ws.Range("M56").Value = -103499.9
$ws.Range("H86").Value = 1490
$ws.Range("I86").Value = 1751.8334
$ws.Range("J86").Value = 966.3333
$ws.Range("K86").Value = 5255.5002
$ws.Range("L86").Value = 2898.9999
$ws.Range("M86").Value = -4069.5002
$ws.Range("N86").Value = -5270.9999
$ws.Range("H89").Value = 1490
$ws.Range("I89").Value = 1751.8334
$ws.Range("J89").Value = 966.3333
$ws.Range("K89").Value = 15766.5006
$ws.Range("L89").Value = 8696.9997
$ws.Range("M89").Value = -9838.500599999999
$ws.Range("N89").Value = -20552.9997
$ws.Range("H122").Value = 1153.2
$ws.Range("I122").Value = 603.0909
$ws.Range("K122").Value = 5427.8181
$ws.Range("M122").Value = -2977.8181
$ws.Range("H133").Value = 4912.0527
$ws.Range("I133").Value = 4873.8
$ws.Range("J133").Value = 4925.7144
$ws.Range("K133").Value = 14621.4
$ws.Range("L133").Value = 14777.1432
$ws.Range("M133").Value = -9561.400000000001
$ws.Range("N133").Value = -24897.1432
$ws.Range("H134").Value = 3547.8684
$ws.Range("I134").Value = 3366.348
$ws.Range("J134").Value = 3826.2
$ws.Range("K134").Value = 10099.044
$ws.Range("L134").Value = 11478.6
$ws.Range("M134").Value = -5029.044
$ws.Range("N134").Value = -21618.6
$ws.Range("H138").Value = 7344.6816
$ws.Range("I138").Value = 13719.875
$ws.Range("J138").Value = 3701.7144
$ws.Range("K138").Value = 41159.625
$ws.Range("L138").Value = 11105.1432
$ws.Range("M138").Value = -36019.625
$ws.Range("N138").Value = -21385.1432
$ws.Range("H140").Value = 1924.9706
$ws.Range("I140").Value = 1688.3
$ws.Range("J140").Value = 3700
$ws.Range("K140").Value = 5064.9
$ws.Range("L140").Value = 11100
$ws.Range("M140").Value = 115.1000000000004
$ws.Range("N140").Value = -21460
$ws.Range("H141").Value = 4151.778
$ws.Range("I141").Value = 3324.4443
$ws.Range("J141").Value = 4979.1113
$ws.Range("K141").Value = 9973.332900000001
$ws.Range("L141").Value = 14937.3339
$ws.Range("M141").Value = -4793.332900000001
$ws.Range("N141").Value = -25297.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 6058.1763
$ws.Range("I80").Value = 18629.666
$ws.Range("J80").Value = 3364.2856
$ws.Range("K80").Value = 18629.666
$ws.Range("L80").Value = 3364.2856
$ws.Range("M80").Value = -17631.666
$ws.Range("N80").Value = -5360.2856
$ws.Range("H83").Value = 6058.1763
$ws.Range("I83").Value = 18629.666
$ws.Range("J83").Value = 3364.2856
$ws.Range("K83").Value = 93148.33
$ws.Range("L83").Value = 16821.428
$ws.Range("M83").Value = -88156.33
$ws.Range("N83").Value = -26805.428
$ws.Range("H132").Value = 42767.32
$ws.Range("I132").Value = 145925.58
$ws.Range("J132").Value = 8381.237999999999
$ws.Range("K132").Value = 437776.74
$ws.Range("L132").Value = 25143.714
$ws.Range("M132").Value = -435246.74
$ws.Range("N132").Value = -30203.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 211140.11
$ws.Range("I55").Value = 333872.16
$ws.Range("J55").Value = 742.2857
$ws.Range("K55").Value = 333872.16
$ws.Range("L55").Value = 742.2857
$ws.Range("M55").Value = -333699.16
$ws.Range("N55").Value = -1088.2857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 67199.336
$ws.Range("I42").Value = 1500
$ws.Range("K42").Value = 1500
$ws.Range("M42").Value = -1122
$ws.Range("H43").Value = 500014750
$ws.Range("I43").Value = 1000000000
$ws.Range("J43").Value = 29500
$ws.Range("K43").Value = 1000000000
$ws.Range("L43").Value = 29500
$ws.Range("M43").Value = -999999851
$ws.Range("N43").Value = -29798
